# Apply edits described by the diff:
# - Clear the numeric values in column B for a number of rows (keep cell styling)
# - Remove cell B30 entirely (value + style)
# - Update B27 value from 21056 to 21086
# - Update sheet view: drop topLeftCell, move selection to B27

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear values but keep styles/formatting on these B-column cells
$rowsToClear = @(3, 7, 9, 10, 12, 13, 14, 15, 16, 18, 19, 20, 21, 22, 28, 31, 32)
foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 2).ClearContents()
}

# B30 loses both its value and its style (the <c> element disappears entirely)
$ws.Cells.Item(30, 2).Clear()

# Update B27 value
$ws.Range("B27").Value = 21086

# Update the active selection/view: select B27 (this also clears any scrolled
# topLeftCell state tied to the previous A28:XFD28 selection)
$ws.Activate()
$ws.Range("B27").Select()
